# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded for the handback status report.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 23:16:38"
$wsZhCn.Range("H2").Value = "2016-03-22 23:17:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 23:16:42"
$wsDeDe.Range("H2").Value = "2016-03-22 23:17:18"
